$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data for 2020-02-13 (homework entry)
$ws.Range("A25").Value = 1581552000

# B25 and C25 must stay as text (date-like string and zero-padded id),
# so force Text format before assigning to avoid Excel's automatic
# date/number reinterpretation.
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "2020-02-13"

$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "03033"

$ws.Range("D25").Value = "SGBHD"
$ws.Range("E25").Value = 0.28
$ws.Range("F25").Value = 0.3
$ws.Range("G25").Value = 0.28
$ws.Range("H25").Value = 0.3
$ws.Range("I25").Value = 120000
